$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (RunNo=1) - update study run details
$ws.Range("E2").Value = "test"
$ws.Range("F2").Value = "10.17.18.109"
$ws.Range("G2").Value = "Test2"
$ws.Range("H2").Value = "2"
$ws.Range("I2").Value = "qualstart-2"

# Row 3 (RunNo=2) - update study run details
$ws.Range("E3").Value = "test"
$ws.Range("F3").Value = "10.17.18.109"
$ws.Range("G3").Value = "Test3"
$ws.Range("H3").Value = "3"
$ws.Range("I3").Value = "qualstart-3"

# Update selection to match the saved state of the workbook (I3 last active cell)
$ws.Range("I3").Select()
